$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("D1").Value = "Total Amount"

# --- Data rows (rows 2-7) ---
$ws.Range("A2").Value = "Yakidoo`n"
$ws.Range("B2").Value = "609`n"
$ws.Range("C2").Value = "3031`n"
$ws.Range("D2").Value = "650,00 `n"

$ws.Range("A3").Value = "Yakidoo`n"
$ws.Range("B3").Value = "609`n"
$ws.Range("C3").Value = "650,00 `n"

$ws.Range("A4").Value = "Jamia`n"
$ws.Range("B4").Value = "97`n"
$ws.Range("C4").Value = "925,00 `n"

$ws.Range("A5").Value = "Mydo`n"
$ws.Range("B5").Value = "570`n"
$ws.Range("C5").Value = "675,00 `n"

$ws.Range("A6").Value = "Tagopia`n"
$ws.Range("B6").Value = "907`n"
$ws.Range("C6").Value = "1 825,00 `n"

$ws.Range("A7").Value = "Fliptune`n"
$ws.Range("B7").Value = "346`n"
$ws.Range("C7").Value = "1 475,00 `n"

Write-Host "values set"
